$wb = $excel.ActiveWorkbook
$ds = $wb.Worksheets.Item("data")

# --- Update time_taken column (F) on the "data" sheet ---
$ds.Range("F2").Value = "2021-10-05 14:21:34.734681"
$ds.Range("F3").Value = "2021-10-05 14:21:34.734691"
$ds.Range("F4").Value = "2021-10-05 14:21:34.734694"
$ds.Range("F5").Value = "2021-10-05 14:21:34.734698"
$ds.Range("F6").Value = "2021-10-05 14:21:34.734701"
$ds.Range("F7").Value = "2021-10-05 14:21:34.734703"
$ds.Range("F8").Value = "2021-10-05 14:21:34.734706"
$ds.Range("F9").Value = "2021-10-05 14:21:34.734709"
$ds.Range("F10").Value = "2021-10-05 14:21:34.734712"
$ds.Range("F11").Value = "2021-10-05 14:21:34.734715"
$ds.Range("F12").Value = "2021-10-05 14:21:34.734718"
$ds.Range("F13").Value = "2021-10-05 14:21:34.734720"
$ds.Range("F14").Value = "2021-10-05 14:21:34.734723"
$ds.Range("F15").Value = "2021-10-05 14:21:34.734726"
$ds.Range("F16").Value = "2021-10-05 14:21:34.734729"
$ds.Range("F17").Value = "2021-10-05 14:21:34.734732"
$ds.Range("F18").Value = "2021-10-05 14:21:34.734735"
$ds.Range("F19").Value = "2021-10-05 14:21:34.734738"
$ds.Range("F20").Value = "2021-10-05 14:21:34.734740"

# --- Add a new "metadata" sheet after "data" ---
$meta = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ds)
$meta.Name = "metadata"

# Header row
$meta.Range("B1").Value = "data_name"
$meta.Range("C1").Value = "data_id"
$meta.Range("D1").Value = "data_version"
$meta.Range("E1").Value = "data_version_created"
$meta.Range("F1").Value = "panel_query_time"
$meta.Range("G1").Value = "panel_get_request"

# Apply the same header/index formatting used on the "data" sheet (the
# bold / centered-top / thin-bordered style already used for its header
# row and row-index column) to the header row and the index cell A2 by
# copying the format from an existing styled cell - this reuses the
# existing style definition instead of creating a near-duplicate one.
$ds.Range("B1").Copy() | Out-Null
$meta.Range("B1:G1").PasteSpecial(-4122) | Out-Null
$meta.Range("A2").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Data row
$meta.Range("A2").Value = 0
$meta.Range("B2").Value = "Mitochondrial disorder with complex V deficiency"
$meta.Range("C2").Value = 538

# Force "1.3" to be stored as text (not a number) while keeping the cell's
# format at the default "Normal" style, same as the source diff.
$meta.Range("D2").NumberFormat = "@"
$meta.Range("D2").Value = "1.3"
$meta.Range("D2").Style = "Normal"

$meta.Range("E2").Value = "2020-02-17T16:12:06.214164Z"
$meta.Range("F2").Value = "2021-10-05 14:21:34.730892"
$meta.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/538/?format=json"

# Keep "data" as the active sheet/tab, matching the source workbook (only
# the sheet list changed; the active tab stayed on "data").
$ds.Activate() | Out-Null
